$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 99 - this shifts the existing rows 99:211 down to 100:212
# and keeps all their data/formatting intact.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly record.
# (Same record as the former row 99 except Fecha, Calidad and Volumen differ.)
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44494
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112023
$ws.Range("G99").Value = "Brócoli"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 500
$ws.Range("L99").Value = 500
$ws.Range("M99").Value = 500
$ws.Range("N99").Value = "$/unidad"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 500
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"
